# Generate Report for Handback
#
# The localization run finished and both target files for the "de-de" and
# "zh-cn" languages have been handed back. This script brings the
# localization-status workbook in line with that: the status text changes
# from "In Translation" to "Handed back: in sync with en-US" everywhere it
# appears, the "Latest Target File" / "Latest Handback File" columns get
# filled in (with a hyperlink on the target-file name, matching the style
# already used for the Source File Name column), and the
# "Latest Handback DateTime" column gets a real timestamp instead of the
# 0001-01-01 placeholder.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96146be47b9337640b0cd40177d8d20de81f582a/e2e/"
$mdName1 = "5496b3a3-cb15-4b22-adf0-bda34f4c4d40.md"
$mdName2 = "fbf4f8ad-7a96-4ba4-87e4-441aef70d8be.md"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (C) for both rows
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Row 2 (5496b3a3...): Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
$wsZh.Range("J2").Value = "5496b3a3-cb15-4b22-adf0-bda34f4c4d40.9a0c6092cef8311764b88743ed3b5fe707803152.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-22 18:25:24"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($baseUrl + $mdName1), "", "", $mdName1)

# Row 3 (fbf4f8ad...): Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
$wsZh.Range("J3").Value = "fbf4f8ad-7a96-4ba4-87e4-441aef70d8be.eb3ffbe5515f3d8f004afcbe97e4a8e16bf533a5.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-22 18:25:24"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($baseUrl + $mdName2), "", "", $mdName2)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column (C) for both rows
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Row 2 (5496b3a3...): Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
$wsDe.Range("J2").Value = "5496b3a3-cb15-4b22-adf0-bda34f4c4d40.9a0c6092cef8311764b88743ed3b5fe707803152.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-22 18:25:32"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($baseUrl + $mdName1), "", "", $mdName1)

# Row 3 (fbf4f8ad...): Latest Target File (I), Latest Handback File (J),
# Latest Handback DateTime (K)
$wsDe.Range("J3").Value = "fbf4f8ad-7a96-4ba4-87e4-441aef70d8be.eb3ffbe5515f3d8f004afcbe97e4a8e16bf533a5.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-22 18:25:32"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($baseUrl + $mdName2), "", "", $mdName2)

# ---------------------------------------------------------------------
# Column widths: the longer status text and the newly-populated file-name
# columns need wider columns so the content isn't truncated (mirrors what
# Excel's AutoFit would do after the content above changed).
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.14
$wsZh.Columns.Item(10).ColumnWidth = 39.14

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.14
$wsDe.Columns.Item(10).ColumnWidth = 39.14
